$d = $word.ActiveDocument

# The new bullet must land right before the existing trailing empty
# "ListParagraph" bullet (currently paragraph 4), so insert a fresh
# paragraph immediately ahead of it; it inherits that paragraph's
# ListParagraph style + numPr (ilvl 0 / numId 1) automatically.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$newParaIndex = $d.Paragraphs.Count - 1
$target = $d.Paragraphs($newParaIndex).Range

# Build the paragraph body as two explicit runs (second run carries a
# leading space, hence xml:space="preserve") via raw OOXML so the two
# runs remain distinct instead of being coalesced into one.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr>' +
       '<w:pStyle w:val="ListParagraph"/>' +
       '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
       '</w:pPr>' +
       '<w:r><w:t>Rename the class. Say Yes to the &quot;You are renaming a file…&quot; dialog.</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> Then make it static.</w:t></w:r>' +
       '</w:p>'

$target.InsertXML($xml)
